$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.422.08"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "1.859.85"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "280.77"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5118"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3520"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.21"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06846"
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.03"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8120"
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07757"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "1.857.72"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.99"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.100"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9982"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008097"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9976"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "26.458.79"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.791"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.211"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.364"
$ws.Range("E25").Value = "  +10.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.44"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.27"
$ws.Range("E28").Value = "  +3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.43"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.372"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.312"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08775"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04892"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.170"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7386"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.865"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.256"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.399"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01859"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5195"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9606"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "115.95"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.267"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.029"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4547"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1360"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.371"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.33"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.504"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").Value = "  +0.61%  "
